$d = $word.ActiveDocument

# --- Step 1: remove the existing _GoBack bookmark (it currently sits right after
#     "...as well as adding ambient noises to the current level tracks") ---
$d.Bookmarks.Item("_GoBack").Delete()

# --- Step 2: locate "put placeholders" freshly (offsets may have shifted after the
#     bookmark deletion, so re-derive them from the live content) ---
$full = $d.Content.Text
$marker = "We ran out of time for the alpha build to put placeholders"
$markerIdx = $full.IndexOf($marker)
$putWord = "put"
$putIdx = $markerIdx + $marker.IndexOf("$putWord placeholders")

# --- Step 3: split the surrounding run into three pieces around "put" by dropping in
#     (and then removing) temporary bookmarks at the word boundaries. Word keeps runs
#     split at a location once a bookmark has forced the split there, even after the
#     bookmark itself is deleted. ---
$d.Bookmarks.Add("zzTempBefore", $d.Range($putIdx, $putIdx)) | Out-Null
$d.Bookmarks.Add("zzTempAfter", $d.Range($putIdx + $putWord.Length, $putIdx + $putWord.Length)) | Out-Null

# --- Step 4: replace "put" with "implement" while the temp bookmarks still protect the
#     new run boundaries from being re-merged with their neighbours ---
$rngPut = $d.Range($putIdx, $putIdx + $putWord.Length)
$rngPut.Text = "implement"

# --- Step 5: remove the temporary bookmarks again; the run split they created persists ---
$d.Bookmarks.Item("zzTempBefore").Delete()
$d.Bookmarks.Item("zzTempAfter").Delete()

# --- Step 6: re-add the _GoBack bookmark right after "implement" (i.e. right where the
#     user's edit just finished, before " placeholders") ---
$afterImplementIdx = $putIdx + "implement".Length
$d.Bookmarks.Add("_GoBack", $d.Range($afterImplementIdx, $afterImplementIdx)) | Out-Null

Write-Host "Done"
